$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").NumberFormat = "mm\-d\-yyyy"
$ws.Range("A3").Value = (Get-Date -Year 2014 -Month 10 -Day 18 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "Add Parse.com and setup the anonymous user. Set the correct device size"

# Row 4
$ws.Range("A4").NumberFormat = "mm\-d\-yyyy"
$ws.Range("A4").Value = (Get-Date -Year 2014 -Month 10 -Day 19 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("B4").Value = 0.5
$ws.Range("C4").Value = "Add include assistant view."

# Row 5
$ws.Range("A5").NumberFormat = "mm\-d\-yyyy"
$ws.Range("A5").Value = (Get-Date -Year 2014 -Month 10 -Day 20 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)

# Row 6
$ws.Range("A6").NumberFormat = "mm\-d\-yyyy"
$ws.Range("A6").Value = (Get-Date -Year 2014 -Month 10 -Day 21 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)

# Row 7
$ws.Range("A7").NumberFormat = "mm\-d\-yyyy"
$ws.Range("A7").Value = (Get-Date -Year 2014 -Month 10 -Day 22 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)

# Adjust column C width and selection (closest achievable value to the
# authored 157.6640625 character-width target given this engine's pixel
# rounding granularity)
$ws.Columns.Item(3).ColumnWidth = 156.833333
$ws.Range("C4").Select()
